$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 486.53333
$ws.Range("I33").Value = 244.55556
$ws.Range("J33").Value = 849.5
$ws.Range("K33").Value = 244.55556
$ws.Range("L33").Value = 849.5
$ws.Range("M33").Value = -15.55556000000001
$ws.Range("N33").Value = -1307.5
$ws.Range("H48").Value = 6000
$ws.Range("I48").Value = 6000
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 18000
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -17708
$ws.Range("H56").Value = 6000
$ws.Range("I56").Value = 6000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 18000
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -17466
$ws.Range("H58").Value = 2056.2354
$ws.Range("I58").Value = 307.125
$ws.Range("J58").Value = 3611
$ws.Range("K58").Value = 921.375
$ws.Range("L58").Value = 10833
$ws.Range("M58").Value = -771.375
$ws.Range("N58").Value = -11133
$ws.Range("H62").Value = 2896.7715
$ws.Range("I62").Value = 2932.9
$ws.Range("J62").Value = 2680
$ws.Range("K62").Value = 2932.9
$ws.Range("L62").Value = 2680
$ws.Range("M62").Value = -2308.9
$ws.Range("N62").Value = -3928
$ws.Range("H65").Value = 2896.7715
$ws.Range("I65").Value = 2932.9
$ws.Range("J65").Value = 2680
$ws.Range("K65").Value = 14664.5
$ws.Range("L65").Value = 13400
$ws.Range("M65").Value = -11544.5
$ws.Range("N65").Value = -19640
$ws.Range("H100").Value = 1407.3077
$ws.Range("I100").Value = 1290.4546
$ws.Range("J100").Value = 2050
$ws.Range("K100").Value = 1290.4546
$ws.Range("L100").Value = 2050
$ws.Range("M100").Value = -749.4546
$ws.Range("N100").Value = -3132
$ws.Range("H111").Value = 337.8
$ws.Range("I111").Value = 337.8
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 1013.4
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = 2053.6
$ws.Range("H138").Value = 30023.541
$ws.Range("I138").Value = 2407.5
$ws.Range("J138").Value = 51064.332
$ws.Range("K138").Value = 7222.5
$ws.Range("L138").Value = 153192.996
$ws.Range("M138").Value = -2082.5
$ws.Range("N138").Value = -163472.996

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 7090.1333
$ws.Range("I5").Value = 8784.75
$ws.Range("J5").Value = 311.66666
$ws.Range("K5").Value = 8784.75
$ws.Range("L5").Value = 311.66666
$ws.Range("M5").Value = -8672.75
$ws.Range("N5").Value = -535.66666
$ws.Range("H32").Value = 42707.19
$ws.Range("I32").Value = 42707.19
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 42707.19
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -42420.19
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H45").Value = 4477
$ws.Range("I45").Value = 3147.8
$ws.Range("J45").Value = 5426.4287
$ws.Range("K45").Value = 3147.8
$ws.Range("L45").Value = 5426.4287
$ws.Range("M45").Value = -2770.8
$ws.Range("N45").Value = -6180.4287
$ws.Range("H88").Value = 6518.091
$ws.Range("I88").Value = 2900
$ws.Range("J88").Value = 9533.166999999999
$ws.Range("K88").Value = 2900
$ws.Range("L88").Value = 9533.166999999999
$ws.Range("M88").Value = -2494
$ws.Range("N88").Value = -10345.167
$ws.Range("H91").Value = 6518.091
$ws.Range("I91").Value = 2900
$ws.Range("J91").Value = 9533.166999999999
$ws.Range("K91").Value = 2900
$ws.Range("L91").Value = 9533.166999999999
$ws.Range("M91").Value = -1496
$ws.Range("N91").Value = -12341.167
$ws.Range("H113").Value = 68000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 68000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 68000
$ws.Range("N113").Value = -76678
$ws.Range("H132").Value = 1425.75
$ws.Range("I132").Value = 1200.3414
$ws.Range("J132").Value = 2746
$ws.Range("K132").Value = 3601.0242
$ws.Range("L132").Value = 8238
$ws.Range("M132").Value = -1071.0242
$ws.Range("N132").Value = -13298

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 7090.1333
$ws.Range("I4").Value = 8784.75
$ws.Range("J4").Value = 311.66666
$ws.Range("K4").Value = 8784.75
$ws.Range("L4").Value = 311.66666
$ws.Range("M4").Value = -8669.75
$ws.Range("N4").Value = -541.66666
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H75").Value = 32852.332
$ws.Range("I75").Value = 28557
$ws.Range("J75").Value = 35000
$ws.Range("K75").Value = 28557
$ws.Range("L75").Value = 35000
$ws.Range("M75").Value = -27621
$ws.Range("N75").Value = -36872
$ws.Range("H78").Value = 32852.332
$ws.Range("I78").Value = 28557
$ws.Range("J78").Value = 35000
$ws.Range("K78").Value = 85671
$ws.Range("L78").Value = 105000
$ws.Range("M78").Value = -80991
$ws.Range("N78").Value = -114360
$ws.Range("H86").Value = 5862.3335
$ws.Range("I86").Value = 1800
$ws.Range("J86").Value = 7893.5
$ws.Range("K86").Value = 1800
$ws.Range("L86").Value = 7893.5
$ws.Range("M86").Value = -677
$ws.Range("N86").Value = -10139.5
$ws.Range("H89").Value = 5862.3335
$ws.Range("I89").Value = 1800
$ws.Range("J89").Value = 7893.5
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 39467.5
$ws.Range("M89").Value = -3384
$ws.Range("N89").Value = -50699.5
$ws.Range("H134").Value = 1784.0454
$ws.Range("I134").Value = 1307.8889
$ws.Range("J134").Value = 3926.75
$ws.Range("K134").Value = 3923.6667
$ws.Range("L134").Value = 11780.25
$ws.Range("M134").Value = -1388.6667
$ws.Range("N134").Value = -16850.25

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 16666.666
$ws.Range("I103").Value = 16666.666
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 16666.666
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -15494.666
$ws.Range("H107").Value = 681.5925999999999
$ws.Range("I107").Value = 575.9
$ws.Range("J107").Value = 983.5714
$ws.Range("K107").Value = 575.9
$ws.Range("L107").Value = 983.5714
$ws.Range("M107").Value = 1344.1
$ws.Range("N107").Value = -4823.5714

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 1183.6666
$ws.Range("I3").Value = 1183.6666
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3550.9998
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -3438.9998
$ws.Range("H63").Value = 6376.96
$ws.Range("I63").Value = 13606
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 40818
$ws.Range("L63").Value = 15000
$ws.Range("M63").Value = -40069
$ws.Range("N63").Value = -16498
$ws.Range("H66").Value = 6376.96
$ws.Range("I66").Value = 13606
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 122454
$ws.Range("L66").Value = 45000
$ws.Range("M66").Value = -118710
$ws.Range("N66").Value = -52488
$ws.Range("H106").Value = 4996.852
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 4996.852
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 14990.556
$ws.Range("N106").Value = -16882.556

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 674.6316
$ws.Range("I2").Value = 879.1667
$ws.Range("J2").Value = 324
$ws.Range("K2").Value = 879.1667
$ws.Range("L2").Value = 324
$ws.Range("M2").Value = -766.1667
$ws.Range("N2").Value = -550
$ws.Range("H102").Value = 17048.908
$ws.Range("I102").Value = 19576.357
$ws.Range("J102").Value = 2895.2
$ws.Range("K102").Value = 19576.357
$ws.Range("L102").Value = 2895.2
$ws.Range("M102").Value = -17954.357
$ws.Range("N102").Value = -6139.2
$ws.Range("H104").Value = 25000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 25000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 25000
$ws.Range("N104").Value = -31988

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1630.75
$ws.Range("I22").Value = 1399.4445
$ws.Range("J22").Value = 1820
$ws.Range("K22").Value = 1399.4445
$ws.Range("L22").Value = 1820
$ws.Range("M22").Value = -1104.4445
$ws.Range("N22").Value = -2410
$ws.Range("H27").Value = 1630.75
$ws.Range("I27").Value = 1399.4445
$ws.Range("J27").Value = 1820
$ws.Range("K27").Value = 1399.4445
$ws.Range("L27").Value = 1820
$ws.Range("M27").Value = -1292.4445
$ws.Range("N27").Value = -2034
$ws.Range("H40").Value = 1838.1428
$ws.Range("I40").Value = 1838.1428
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1838.1428
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1702.1428
$ws.Range("H46").Value = 2814.9583
$ws.Range("I46").Value = 802.5
$ws.Range("J46").Value = 3217.45
$ws.Range("K46").Value = 802.5
$ws.Range("L46").Value = 3217.45
$ws.Range("M46").Value = -614.5
$ws.Range("N46").Value = -3593.45
$ws.Range("H55").Value = 1130.6666
$ws.Range("I55").Value = 115.833336
$ws.Range("J55").Value = 2145.5
$ws.Range("K55").Value = 115.833336
$ws.Range("L55").Value = 2145.5
$ws.Range("M55").Value = 57.166664
$ws.Range("N55").Value = -2491.5
$ws.Range("H93").Value = 1413.4231
$ws.Range("I93").Value = 1054.6316
$ws.Range("J93").Value = 2387.2856
$ws.Range("K93").Value = 1054.6316
$ws.Range("L93").Value = 2387.2856
$ws.Range("M93").Value = 193.3684000000001
$ws.Range("N93").Value = -4883.2856

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 30000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 30000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30766
$ws.Range("H85").Value = 30000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 30000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32652
$ws.Range("H122").Value = 108238.61
$ws.Range("I122").Value = 137005.83
$ws.Range("J122").Value = 4676.6
$ws.Range("K122").Value = 411017.49
$ws.Range("L122").Value = 14029.8
$ws.Range("M122").Value = -408567.49
$ws.Range("N122").Value = -18929.8
